$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# "Corinne's data" columns got renamed from the old Vulture ID (CK151400) to the
# new bird id (#109018542) - update the fisher-formula header labels accordingly.
$ws.Range("K1").Value = "#109018542 x_"
$ws.Range("L1").Value = "#109018542 y_"

# Scroll the sheet over (so column D is at the left edge) and leave the
# selection on K25, matching where the author was working.
$excel.ActiveWindow.ScrollColumn = 4
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("K25").Select()
